# Update TPM-derived NATMI ligand-receptor metrics (L1cam-Cd9) with the
# newly recomputed values (commit: "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("M2").Value = 19.515399
$ws.Range("N2").Value = 58.546197
$ws.Range("O2").Value = 0.1046357846766865
$ws.Range("P2").Value = 0.1046357846766865
$ws.Range("Q2").Value = 151.069269605571
$ws.Range("R2").Value = 1359.623426450139
$ws.Range("S2").Value = 0.05159053691576445
$ws.Range("T2").Value = 0.05159053691576446
$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("M3").Value = 133.0753813333333
$ws.Range("N3").Value = 399.226144
$ws.Range("O3").Value = 0.7135107484588257
$ws.Range("P3").Value = 0.7135107484588257
$ws.Range("Q3").Value = 1030.140386087392
$ws.Range("R3").Value = 9271.263474786527
$ws.Range("S3").Value = 0.3517955422411176
$ws.Range("T3").Value = 0.3517955422411176
$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.1818534668644878
$ws.Range("P4").Value = 0.1818534668644878
$ws.Range("Q4").Value = 262.553298561731
$ws.Range("R4").Value = 2362.979687055579
$ws.Range("S4").Value = 0.0896626142243902
$ws.Range("T4").Value = 0.0896626142243902
$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("M5").Value = 19.515399
$ws.Range("N5").Value = 58.546197
$ws.Range("O5").Value = 0.1046357846766865
$ws.Range("P5").Value = 0.1046357846766865
$ws.Range("Q5").Value = 5.944364514868
$ws.Range("R5").Value = 53.499280633812
$ws.Range("S5").Value = 0.002030015487238104
$ws.Range("T5").Value = 0.002030015487238105
$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("M6").Value = 133.0753813333333
$ws.Range("N6").Value = 399.226144
$ws.Range("O6").Value = 0.7135107484588257
$ws.Range("P6").Value = 0.7135107484588257
$ws.Range("R6").Value = 364.811253482624
$ws.Range("S6").Value = 0.01384266266227932
$ws.Range("T6").Value = 0.01384266266227932
$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.1818534668644878
$ws.Range("P7").Value = 0.1818534668644878
$ws.Range("S7").Value = 0.003528098492150974
$ws.Range("T7").Value = 0.003528098492150974
$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("M8").Value = 19.515399
$ws.Range("N8").Value = 58.546197
$ws.Range("O8").Value = 0.1046357846766865
$ws.Range("P8").Value = 0.1046357846766865
$ws.Range("Q8").Value = 149.384641817694
$ws.Range("R8").Value = 1344.461776359246
$ws.Range("S8").Value = 0.05101523227368399
$ws.Range("T8").Value = 0.05101523227368399
$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("M9").Value = 133.0753813333333
$ws.Range("N9").Value = 399.226144
$ws.Range("O9").Value = 0.7135107484588257
$ws.Range("P9").Value = 0.7135107484588257
$ws.Range("S9").Value = 0.3478725435554287
$ws.Range("T9").Value = 0.3478725435554287
$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.1818534668644878
$ws.Range("P10").Value = 0.1818534668644878
$ws.Range("S10").Value = 0.08866275414794662
$ws.Range("T10").Value = 0.08866275414794662
